# New Submission Synced: 2026-02-08 19:09:47
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# Fix up C5, which was stored as text ("32") but should be a real number
$ws.Range("C5").Value = 32

# Append the new submission row
$ws.Range("A6").Value = "2026-02-08 19:09:47"
$ws.Range("B6").Value = "Iyudigal Fali"
# Force "4" to be stored as text (admission numbers are text in this sheet),
# then restore the default style so no extra formatting sticks to the cell.
$ws.Range("C6").Value = "'4"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 9
